# edit.ps1 - applies "curso avanzado de laravel 11" commit changes
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Actualmente trabajo..." paragraph - trim the tech stack list
#   ", Node.js, React.js, Vue.js"  ->  ", Vue.js"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    ", Node.js, React.js, Vue.js", $true, $false, $false, $false, $false,
    $true, 1, $false, ", Vue.js", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: "Informacion general" paragraph - update tech stack mention and
# split it out into new runs for "Vue.js" and "MySQL"
#   "...Laravel, Node.js, React.js y AWS. Gran sentido de..."
#   -> "...Laravel, Vue.js y MySQL. Gran sentido de..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Laravel, Node.js, React.js y AWS. Gran sentido de ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Laravel, Vue.js y MySQL. Gran sentido de ", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: append the new "Carta de motivacion" section at the end of the
# document, right before the trailing empty paragraph.
# ---------------------------------------------------------------------------
function Insert-ParaBefore([string]$text) {
    $count = $d.Paragraphs.Count
    $target = $d.Paragraphs.Item($count)
    $pos = $target.Range.Start
    $r = $d.Range($pos, $pos)
    $r.InsertParagraphBefore() | Out-Null
    $newCount = $d.Paragraphs.Count
    $p = $d.Paragraphs.Item($newCount - 1)
    if ($text) {
        $p.Range.InsertAfter($text) | Out-Null
    }
    return $p
}

# empty paragraph
Insert-ParaBefore("") | Out-Null

# bold heading paragraph
$headingPara = Insert-ParaBefore("Carta de motivación:")
$headingPara.Range.Font.Bold = $true

Insert-ParaBefore("A lo largo de mi trayectoria profesional, he cultivado una pasión por la programación y las tecnologías web, lo cual me ha permitido adquirir una sólida experiencia en el uso de herramientas y lenguajes como Laravel, Node.js, React.js, Vue.js, SQL y AWS.") | Out-Null

Insert-ParaBefore("Mi formación como Ingeniero Químico de la Universidad Central de Venezuela me dotó de una mentalidad analítica y una sólida base en programación y bases de datos, que he complementado con una amplia capacitación en tecnologías web a través de cursos acreditados por instituciones reconocidas. Este constante afán por el aprendizaje me ha permitido mantenerme a la vanguardia en un campo en constante evolución.") | Out-Null

Insert-ParaBefore("Durante más de 15 años en la industria petrolera, demostré mi capacidad para automatizar documentos y procesos complejos mediante el uso de VBA, VB6, Delphi, Access, Excel, MySQL y SQL Server. Esta experiencia me ha dotado de una habilidad excepcional para identificar oportunidades de optimización y desarrollar soluciones innovadoras que generan valor agregado para las organizaciones.") | Out-Null

Insert-ParaBefore("Más recientemente, he tenido la oportunidad de aplicar mis conocimientos en el desarrollo web, trabajando desde agosto de 2019 en proyectos que involucran tecnologías como Laravel, Node.js, React.js, Vue.js, SQL y AWS. Esta experiencia me ha permitido adquirir una comprensión profunda de los desafíos y las mejores prácticas en el desarrollo de aplicaciones web modernas y escalables.") | Out-Null

Insert-ParaBefore("Además de mi experiencia técnica, me distingue una gran capacidad para trabajar bajo presión, un alto sentido de responsabilidad y una vocación de servicio inquebrantable. Soy una persona dinámica, capaz de aportar ideas innovadoras y tomar la iniciativa en la resolución de problemas, siempre con el objetivo de generar valor para la organización.") | Out-Null

Insert-ParaBefore("Estoy convencido de que mis habilidades técnicas, mi versatilidad y mi compromiso con el aprendizaje continuo me convierten en un candidato excepcional para formar parte de su equipo. Anhelo tener la oportunidad de contribuir con mi experiencia y pasión al éxito de su empresa.") | Out-Null

Insert-ParaBefore("Gracias por su consideración. Espero con gran ilusión la oportunidad de discutir personalmente cómo puedo aportar valor a su organización.") | Out-Null
